# Update the example BOM data on Sheet1.
# Part numbers were renamed/renumbered; quantities (column B) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A3/A4 first so the shared-string table is rebuilt in the same order
# as the authoritative edit (WH-01 already existed, TR-01 is added next,
# then the new SK10xx-01 values).
$ws.Range("A3").Value = "WH-01"
$ws.Range("A4").Value = "TR-01"
$ws.Range("A2").Value = "SK1002-01"
$ws.Range("A5").Value = "SK1005-01"
$ws.Range("A6").Value = "SK1007-01"

# Column A needs to be a little wider to fit the new part numbers.
$ws.Columns("A").ColumnWidth = 14.57

# The active selection moved to B4.
$ws.Range("B4").Select() | Out-Null
